$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinate values in row 2 to whole numbers
$ws.Range("Q2").Value = 798924
$ws.Range("R2").Value = 7235332

# Remove the Starttid (Z2) and Sluttid (AB2) time values for row 2
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Round the easting/northing coordinate values in row 3 to whole numbers
$ws.Range("Q3").Value = 798928
$ws.Range("R3").Value = 7235310

# Remove the Starttid (Z3) and Sluttid (AB3) time values for row 3
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
